$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 443, shifting the existing rows 443:463 down to 444:464
$ws.Rows.Item(443).Insert()

# Populate the newly inserted row 443 with the new weekly data point
$ws.Cells.Item(443, 1).Value2  = 4
$ws.Cells.Item(443, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(443, 3).Value2  = "Los Lagos"
$ws.Cells.Item(443, 4).Value2  = 45041
$ws.Cells.Item(443, 5).Value2  = 10
$ws.Cells.Item(443, 6).Value2  = 100112045
$ws.Cells.Item(443, 7).Value2  = "Zapallo"
$ws.Cells.Item(443, 8).Value2  = "Paine"
$ws.Cells.Item(443, 9).Value2  = "1a (cosecha)"
$ws.Cells.Item(443, 10).Value2 = 1200
$ws.Cells.Item(443, 11).Value2 = 550
$ws.Cells.Item(443, 12).Value2 = 600
$ws.Cells.Item(443, 13).Value2 = 575
$ws.Cells.Item(443, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(443, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(443, 16).Value2 = 575
$ws.Cells.Item(443, 17).Value2 = 1
$ws.Cells.Item(443, 18).Value2 = "Hortaliza"
